# Auto-generated Excel COM-interop script
# Applies numeric updates to the Asura_Profits price/profit sheets
# (currentAveragePrice / *NQ / *HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# as refreshed by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 11977563
$ws.Range("J70").Value = 1362.625
$ws.Range("L70").Value = 4087.875
$ws.Range("N70").Value = -4627.875

$ws.Range("H73").Value = 11977563
$ws.Range("J73").Value = 1362.625
$ws.Range("L73").Value = 4087.875
$ws.Range("N73").Value = -5959.875

$ws.Range("H75").Value = 32666.666
$ws.Range("J75").Value = 34000
$ws.Range("L75").Value = 34000
$ws.Range("N75").Value = -35872

$ws.Range("H78").Value = 32666.666
$ws.Range("J78").Value = 34000
$ws.Range("L78").Value = 102000
$ws.Range("N78").Value = -111360

$ws.Range("H98").Value = 7636.2354
$ws.Range("I98").Value = 5070.385
$ws.Range("J98").Value = 15975.25
$ws.Range("K98").Value = 5070.385
$ws.Range("L98").Value = 15975.25
$ws.Range("M98").Value = -3572.385
$ws.Range("N98").Value = -18971.25

$ws.Range("H112").Value = 8129.7646
$ws.Range("J112").Value = 11142.167
$ws.Range("L112").Value = 33426.501
$ws.Range("N112").Value = -35642.501

$ws.Range("H122").Value = 7636.2354
$ws.Range("I122").Value = 5070.385
$ws.Range("J122").Value = 15975.25
$ws.Range("K122").Value = 15211.155
$ws.Range("L122").Value = 47925.75
$ws.Range("M122").Value = -12761.155
$ws.Range("N122").Value = -52825.75

$ws.Range("H132").Value = 1247.3438
$ws.Range("I132").Value = 1202.4674
$ws.Range("J132").Value = 2279.5
$ws.Range("K132").Value = 3607.4022
$ws.Range("L132").Value = 6838.5
$ws.Range("M132").Value = -1077.4022
$ws.Range("N132").Value = -11898.5

$ws.Range("H138").Value = 2045533.5
$ws.Range("I138").Value = 4002812.2
$ws.Range("J138").Value = 6701.4375
$ws.Range("K138").Value = 12008436.6
$ws.Range("L138").Value = 20104.3125
$ws.Range("M138").Value = -12003296.6
$ws.Range("N138").Value = -30384.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41554
$ws.Range("I32").Value = 37157.555
$ws.Range("K32").Value = 37157.555
$ws.Range("M32").Value = -36870.555

$ws.Range("H61").Value = 2770.4348
$ws.Range("I61").Value = 2393.647
$ws.Range("J61").Value = 3838
$ws.Range("K61").Value = 2393.647
$ws.Range("L61").Value = 3838
$ws.Range("M61").Value = -2181.647
$ws.Range("N61").Value = -4262

$ws.Range("H74").Value = 1304.6666
$ws.Range("I74").Value = 1278.9667
$ws.Range("J74").Value = 1433.1666
$ws.Range("K74").Value = 1278.9667
$ws.Range("L74").Value = 1433.1666
$ws.Range("M74").Value = -404.9666999999999
$ws.Range("N74").Value = -3181.1666

$ws.Range("H77").Value = 1304.6666
$ws.Range("I77").Value = 1278.9667
$ws.Range("J77").Value = 1433.1666
$ws.Range("K77").Value = 6394.8335
$ws.Range("L77").Value = 7165.833000000001
$ws.Range("M77").Value = -2026.8335
$ws.Range("N77").Value = -15901.833

$ws.Range("H122").Value = 3542.4
$ws.Range("I122").Value = 3178
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9534
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -7084
$ws.Range("N122").Value = -19900

$ws.Range("H136").Value = 2770.4348
$ws.Range("I136").Value = 2393.647
$ws.Range("J136").Value = 3838
$ws.Range("K136").Value = 7180.941
$ws.Range("L136").Value = 11514
$ws.Range("M136").Value = -4630.941
$ws.Range("N136").Value = -16614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 98750
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 98750
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 98750
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -100190

$ws.Range("H136").Value = 98750
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 98750
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 98750
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -108950

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5495
$ws.Range("I4").Value = 990
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 990
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = -878
$ws.Range("N4").Value = -10224

$ws.Range("H31").Value = 3545.9167
$ws.Range("I31").Value = 2409.625
$ws.Range("J31").Value = 5818.5
$ws.Range("K31").Value = 2409.625
$ws.Range("L31").Value = 5818.5
$ws.Range("M31").Value = -2114.625
$ws.Range("N31").Value = -6408.5

$ws.Range("H34").Value = 3545.9167
$ws.Range("I34").Value = 2409.625
$ws.Range("J34").Value = 5818.5
$ws.Range("K34").Value = 2409.625
$ws.Range("L34").Value = 5818.5
$ws.Range("M34").Value = -2207.625
$ws.Range("N34").Value = -6222.5

$ws.Range("H42").Value = 36666.668
$ws.Range("I42").Value = 32500
$ws.Range("K42").Value = 32500
$ws.Range("M42").Value = -31907

$ws.Range("H69").Value = 9600.5
$ws.Range("I69").Value = 4000
$ws.Range("J69").Value = 15201
$ws.Range("K69").Value = 4000
$ws.Range("L69").Value = 15201
$ws.Range("M69").Value = -3251
$ws.Range("N69").Value = -16699

$ws.Range("H72").Value = 9600.5
$ws.Range("I72").Value = 4000
$ws.Range("J72").Value = 15201
$ws.Range("K72").Value = 12000
$ws.Range("L72").Value = 45603
$ws.Range("M72").Value = -8256
$ws.Range("N72").Value = -53091

$ws.Range("H122").Value = 2506
$ws.Range("I122").Value = 2755.5
$ws.Range("K122").Value = 8266.5
$ws.Range("M122").Value = -5816.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 178398
$ws.Range("I4").Value = 250097
$ws.Range("J4").Value = 35000
$ws.Range("K4").Value = 750291
$ws.Range("L4").Value = 105000
$ws.Range("M4").Value = -750179
$ws.Range("N4").Value = -105224

$ws.Range("H10").Value = 2164.1428
$ws.Range("I10").Value = 62.5
$ws.Range("K10").Value = 187.5
$ws.Range("M10").Value = -48.5

$ws.Range("H17").Value = 724.8333
$ws.Range("I17").Value = 449
$ws.Range("J17").Value = 780
$ws.Range("K17").Value = 1347
$ws.Range("L17").Value = 2340
$ws.Range("M17").Value = -1178
$ws.Range("N17").Value = -2678

$ws.Range("H34").Value = 929.4375
$ws.Range("J34").Value = 1057.1482
$ws.Range("L34").Value = 3171.4446
$ws.Range("N34").Value = -3339.4446

$ws.Range("H39").Value = 2698.7368
$ws.Range("J39").Value = 2933.8823
$ws.Range("L39").Value = 8801.6469
$ws.Range("N39").Value = -9389.6469

$ws.Range("H55").Value = 3261.3845
$ws.Range("J55").Value = 3261.3845
$ws.Range("L55").Value = 9784.1535
$ws.Range("N55").Value = -10138.1535

$ws.Range("H62").Value = 1750
$ws.Range("J62").Value = 2500
$ws.Range("L62").Value = 7500
$ws.Range("N62").Value = -8872

$ws.Range("H63").Value = 177818.25
$ws.Range("I63").Value = 234979.89
$ws.Range("J63").Value = 6333.3335
$ws.Range("K63").Value = 704939.67
$ws.Range("L63").Value = 19000.0005
$ws.Range("M63").Value = -704190.67
$ws.Range("N63").Value = -20498.0005

$ws.Range("H65").Value = 1750
$ws.Range("J65").Value = 2500
$ws.Range("L65").Value = 22500
$ws.Range("N65").Value = -29364

$ws.Range("H66").Value = 177818.25
$ws.Range("I66").Value = 234979.89
$ws.Range("J66").Value = 6333.3335
$ws.Range("K66").Value = 2114819.01
$ws.Range("L66").Value = 57000.0015
$ws.Range("M66").Value = -2111075.01
$ws.Range("N66").Value = -64488.0015

$ws.Range("H69").Value = 873.8182
$ws.Range("I69").Value = 506
$ws.Range("J69").Value = 955.55554
$ws.Range("K69").Value = 1518
$ws.Range("L69").Value = 2866.66662
$ws.Range("M69").Value = -707
$ws.Range("N69").Value = -4488.66662

$ws.Range("H70").Value = 19000
$ws.Range("I70").Value = 19000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 57000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -56685
$ws.Range("N70").ClearContents()

$ws.Range("H72").Value = 873.8182
$ws.Range("I72").Value = 506
$ws.Range("J72").Value = 955.55554
$ws.Range("K72").Value = 4554
$ws.Range("L72").Value = 8599.99986
$ws.Range("M72").Value = -498
$ws.Range("N72").Value = -16711.99986

$ws.Range("H73").Value = 19000
$ws.Range("I73").Value = 19000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 57000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -55908
$ws.Range("N73").ClearContents()

$ws.Range("H131").Value = 869.1111
$ws.Range("I131").Value = 687.5
$ws.Range("J131").Value = 876.7578999999999
$ws.Range("K131").Value = 2062.5
$ws.Range("L131").Value = 2630.2737
$ws.Range("M131").Value = 2977.5
$ws.Range("N131").Value = -12710.2737

$ws.Range("H132").Value = 1851.1177
$ws.Range("I132").Value = 1451.8
$ws.Range("J132").Value = 2017.5
$ws.Range("K132").Value = 13066.2
$ws.Range("L132").Value = 18157.5
$ws.Range("M132").Value = -10536.2
$ws.Range("N132").Value = -23217.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 6320.4
$ws.Range("J36").Value = 7528.909
$ws.Range("L36").Value = 7528.909
$ws.Range("N36").Value = -8498.909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9995
$ws.Range("J2").Value = 9995
$ws.Range("L2").Value = 9995
$ws.Range("N2").Value = -10219

$ws.Range("H122").Value = 13976818
$ws.Range("I122").Value = 16442021
$ws.Range("J122").Value = 7335
$ws.Range("K122").Value = 49326063
$ws.Range("L122").Value = 22005
$ws.Range("M122").Value = -49323613
$ws.Range("N122").Value = -26905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -31240

$ws.Range("H106").Value = 33000
$ws.Range("J106").Value = 33000
$ws.Range("L106").Value = 33000
$ws.Range("N106").Value = -35524
